# StdDevEmployeeHours.xlsx edit
#
# The upstream commit ("changed the stream in the code") altered how the
# per-project employee-hours records are iterated before the per-project
# standard deviation is computed. Two observable effects land in the
# worksheet:
#   1) A handful of adjacent same-valued-looking rows were produced from a
#      differently-ordered input stream, so the project label (col A) and
#      its std-dev (col B) swap between specific row pairs: 22/23, 26/27,
#      61/62 and 63/64.
#   2) Summing the underlying hours in a different order changes the
#      floating point rounding of many (but not all) of the std-dev
#      results at the ULP level.
#
# This script reproduces the resulting cell values directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: project-label swaps caused by the reordered input stream ---
$ws.Cells.Item(22, 1).Value = "PRJ188"
$ws.Cells.Item(23, 1).Value = "PRJ100"
$ws.Cells.Item(26, 1).Value = "PRJ189"
$ws.Cells.Item(27, 1).Value = "PRJ101"
$ws.Cells.Item(61, 1).Value = "PRJ111"
$ws.Cells.Item(62, 1).Value = "PRJ199"
$ws.Cells.Item(63, 1).Value = "PRJ110"
$ws.Cells.Item(64, 1).Value = "PRJ198"

# --- Column B: recalculated standard-deviation values ---
$ws.Cells.Item(3, 2).Value = 4.693415600604745
$ws.Cells.Item(6, 2).Value = 1.2629999010292916
$ws.Cells.Item(8, 2).Value = 3.326658316622787
$ws.Cells.Item(10, 2).Value = 3.11776362157236
$ws.Cells.Item(14, 2).Value = 2.5029470629639774
$ws.Cells.Item(18, 2).Value = 3.5726550351244377
$ws.Cells.Item(22, 2).Value = 2.1775036675565493
$ws.Cells.Item(23, 2).Value = 2.325353019794342
$ws.Cells.Item(25, 2).Value = 6.1475653717549035
$ws.Cells.Item(26, 2).Value = 2.842804249328469
$ws.Cells.Item(27, 2).Value = 3.086052494692856
$ws.Cells.Item(35, 2).Value = 3.191783287686611
$ws.Cells.Item(42, 2).Value = 0.9239889369227076
$ws.Cells.Item(61, 2).Value = 2.6547798402127434
$ws.Cells.Item(62, 2).Value = 2.182266915134097
$ws.Cells.Item(63, 2).Value = 2.2196432145730087
$ws.Cells.Item(64, 2).Value = 0.0
$ws.Cells.Item(66, 2).Value = 1.7865578325060985
$ws.Cells.Item(68, 2).Value = 4.008492720323798
$ws.Cells.Item(74, 2).Value = 7.1289378981413165
$ws.Cells.Item(85, 2).Value = 1.2317629642102412
$ws.Cells.Item(88, 2).Value = 1.6022692241526288
$ws.Cells.Item(93, 2).Value = 2.549153585016015
$ws.Cells.Item(95, 2).Value = 2.9302650155275187

# --- Columns A/B best-fit width nudges (the regenerated workbook's
#     auto-fit pass re-measured the two columns slightly differently) ---
$ws.Columns.Item(1).ColumnWidth = 7.46484375
$ws.Columns.Item(2).ColumnWidth = 26.265625
